# The new weekly "Jengibre" price observation needs to be inserted as the
# first data row of this block (row 94), pushing the existing rows 94:219
# down by one (to 95:220) — i.e. a normal row insert followed by filling
# in the new row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 94; everything below (old rows 94-219) shifts down
# to 95-220, carrying its formatting (incl. the date-formatted column D).
$ws.Rows("94").Insert()

# Populate the newly inserted row 94 with the latest observation.
$ws.Range("A94").Value = 10
$ws.Range("B94").Value = "Vega Modelo de Temuco"
$ws.Range("C94").Value = "La Araucanía"
$ws.Range("D94").Value = 44897
$ws.Range("E94").Value = 9
$ws.Range("F94").Value = 100114007
$ws.Range("G94").Value = "Jengibre"
$ws.Range("H94").Value = "Sin especificar"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 80
$ws.Range("K94").Value = 20000
$ws.Range("L94").Value = 20000
$ws.Range("M94").Value = 20000
$ws.Range("N94").Value = "`$/caja 13 kilos"
$ws.Range("O94").Value = "Perú"
$ws.Range("P94").Value = 1538
$ws.Range("Q94").Value = 13
$ws.Range("R94").Value = "Hortaliza"
